$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 24, shifting the existing data rows (24-141) down to (25-142)
$ws.Rows.Item(24).Insert()

# Populate the newly inserted row 24 with the new price record
$ws.Cells.Item(24, 1).Value = 4
$ws.Cells.Item(24, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(24, 3).Value = "Los Lagos"
$ws.Cells.Item(24, 4).Value = 44831
$ws.Cells.Item(24, 5).Value = 10
$ws.Cells.Item(24, 6).Value = 100112052
$ws.Cells.Item(24, 7).Value = "Albahaca"
$ws.Cells.Item(24, 8).Value = "Sin especificar"
$ws.Cells.Item(24, 9).Value = "Primera"
$ws.Cells.Item(24, 10).Value = 80
$ws.Cells.Item(24, 11).Value = 6000
$ws.Cells.Item(24, 12).Value = 6000
$ws.Cells.Item(24, 13).Value = 6000
$ws.Cells.Item(24, 14).Value = "`$/docena de matas"
$ws.Cells.Item(24, 15).Value = "Región Metropolitana"
$ws.Cells.Item(24, 16).Value = 1000
$ws.Cells.Item(24, 17).Value = 6
$ws.Cells.Item(24, 18).Value = "Hortaliza"
